$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.112.64"
$ws.Range("E2").Value = "  +0.11%  "
# Row 3
$ws.Range("D3").Value = "1.882.87"
$ws.Range("E3").Value = "  -0.96%  "
# Row 4
$ws.Range("E4").Value = "  +0.16%  "
# Row 5
$ws.Range("D5").Value = "'313.55"
$ws.Range("E5").Value = "  +0.37%  "
# Row 7
$ws.Range("D7").Value = "'0.5082"
$ws.Range("E7").Value = "  +0.77%  "
# Row 8
$ws.Range("D8").Value = "'0.3869"
$ws.Range("E8").Value = "  -1.41%  "
# Row 9
$ws.Range("D9").Value = "'0.09027"
$ws.Range("E9").Value = "  -3.41%  "
# Row 10
$ws.Range("D10").Value = "'1.126"
$ws.Range("E10").Value = "  -0.88%  "
# Row 11
$ws.Range("D11").Value = "'41.68"
$ws.Range("E11").Value = "  -0.15%  "
# Row 12
$ws.Range("D12").Value = "'6.365"
$ws.Range("E12").Value = "  -0.14%  "
# Row 13
$ws.Range("D13").Value = "'20.81"
$ws.Range("E13").Value = "  +0.17%  "
# Row 14
$ws.Range("D14").Value = "1.870.98"
$ws.Range("E14").Value = "  -1.08%  "
# Row 15
$ws.Range("D15").Value = "'7.263"
$ws.Range("E15").Value = "  -0.68%  "
# Row 16
$ws.Range("E16").Value = "  +0.18%  "
# Row 17
$ws.Range("E17").Value = "  -0.56%  "
# Row 18
$ws.Range("D18").Value = "'91.42"
$ws.Range("E18").Value = "  -1.20%  "
# Row 19
$ws.Range("D19").Value = "'0.06629"
$ws.Range("E19").Value = "  +0.84%  "
# Row 20
$ws.Range("D20").Value = "'18.23"
$ws.Range("E20").Value = "  +2.28%  "
# Row 21
$ws.Range("E21").Value = "  +0.19%  "
# Row 22
$ws.Range("D22").Value = "'6.127"
# Row 23
$ws.Range("D23").Value = "28.147.96"
$ws.Range("E23").Value = "  +0.05%  "
# Row 24
$ws.Range("D24").Value = "'11.44"
$ws.Range("E24").Value = "  +0.63%  "
# Row 25
$ws.Range("E25").Value = "  -2.22%  "
# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'3.387"
$ws.Range("E26").Value = "  +0.09%  "
# Row 27
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.100.65"
$ws.Range("E27").Value = "  -0.21%  "
# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.548"
$ws.Range("E28").Value = "  -3.15%  "
# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.81"
$ws.Range("E29").Value = "  -0.36%  "
# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'157.19"
$ws.Range("E30").Value = "  -0.02%  "
# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'127.18"
$ws.Range("E31").Value = "  +0.02%  "
# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1059"
$ws.Range("E32").Value = "  -0.61%  "
# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.066"
$ws.Range("E33").Value = "  -1.73%  "
# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.625"
$ws.Range("E34").Value = "  +0.06%  "
# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.602"
$ws.Range("E35").Value = "  -0.48%  "
# Row 36
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.591"
$ws.Range("E36").Value = "  -0.53%  "
# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06633"
$ws.Range("E37").Value = "  -0.17%  "
# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02410"
$ws.Range("E38").Value = "  -0.50%  "
# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2189"
$ws.Range("E39").Value = "  +0.63%  "
# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.286"
$ws.Range("E40").Value = "  +1.09%  "
# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.215"
$ws.Range("E41").Value = "  -1.11%  "
# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6416"
$ws.Range("E42").Value = "  +0.47%  "
# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'11.52"
$ws.Range("E43").Value = "  +0.77%  "
# Row 44
$ws.Range("B44").Value = "InternetComputer(DFINITY)"
$ws.Range("C44").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D44").Value = "'4.926"
$ws.Range("E44").Value = "  -1.42%  "
# Row 45
$ws.Range("D45").Value = "'0.6054"
$ws.Range("E45").Value = "  +0.93%  "
# Row 46
$ws.Range("D46").Value = "'13.19"
$ws.Range("E46").Value = "  -1.09%  "
# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.669"
$ws.Range("E47").Value = "  -1.25%  "
# Row 48
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.275"
$ws.Range("E48").Value = "  +0.11%  "
# Row 49
$ws.Range("D49").Value = "'1.247"
$ws.Range("E49").Value = "  +5.83%  "
# Row 50
$ws.Range("D50").Value = "'2.004"
$ws.Range("E50").Value = "  -1.13%  "
# Row 51
$ws.Range("D51").Value = "'121.46"
$ws.Range("E51").Value = "  -0.58%  "
